$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each entry represents one cell whose text content changes in this update
# (updated Price / Volume(1h) figures for the crypto list).
$updates = @(
    @{ Cell = 'D2'; Value = '22.014.66' },
    @{ Cell = 'E2'; Value = '  -1.10%  ' },
    @{ Cell = 'D3'; Value = '1.552.49' },
    @{ Cell = 'E3'; Value = '  -0.31%  ' },
    @{ Cell = 'D4'; Value = '1.003' },
    @{ Cell = 'E4'; Value = '  +0.10%  ' },
    @{ Cell = 'E5'; Value = '  +0.06%  ' },
    @{ Cell = 'D6'; Value = '290.17' },
    @{ Cell = 'E6'; Value = '  +0.42%  ' },
    @{ Cell = 'D7'; Value = '0.3926' },
    @{ Cell = 'E7'; Value = '  +3.30%  ' },
    @{ Cell = 'D8'; Value = '0.3207' },
    @{ Cell = 'E8'; Value = '  -2.34%  ' },
    @{ Cell = 'D9'; Value = '43.43' },
    @{ Cell = 'E9'; Value = '  -2.90%  ' },
    @{ Cell = 'D10'; Value = '0.07229' },
    @{ Cell = 'E10'; Value = '  -2.10%  ' },
    @{ Cell = 'D11'; Value = '1.068' },
    @{ Cell = 'E11'; Value = '  -6.22%  ' },
    @{ Cell = 'E12'; Value = '  +0.10%  ' },
    @{ Cell = 'D13'; Value = '5.650' },
    @{ Cell = 'E13'; Value = '  -3.60%  ' },
    @{ Cell = 'D14'; Value = '18.62' },
    @{ Cell = 'E14'; Value = '  -8.30%  ' },
    @{ Cell = 'D15'; Value = '0.00001125' },
    @{ Cell = 'E15'; Value = '  +4.46%  ' },
    @{ Cell = 'D16'; Value = '6.597' },
    @{ Cell = 'E16'; Value = '  -2.53%  ' },
    @{ Cell = 'D17'; Value = '1.551.68' },
    @{ Cell = 'E17'; Value = '  -0.19%  ' },
    @{ Cell = 'D18'; Value = '0.06583' },
    @{ Cell = 'E18'; Value = '  -1.19%  ' },
    @{ Cell = 'D19'; Value = '83.24' },
    @{ Cell = 'E19'; Value = '  -3.75%  ' },
    @{ Cell = 'E20'; Value = '  -0.21%  ' },
    @{ Cell = 'D21'; Value = '6.267' },
    @{ Cell = 'E21'; Value = '  -2.71%  ' },
    @{ Cell = 'D22'; Value = '15.40' },
    @{ Cell = 'E22'; Value = '  -5.01%  ' },
    @{ Cell = 'D23'; Value = '11.22' },
    @{ Cell = 'E23'; Value = '  -4.55%  ' },
    @{ Cell = 'D24'; Value = '22.029.56' },
    @{ Cell = 'E24'; Value = '  -1.02%  ' },
    @{ Cell = 'D25'; Value = '2.379' },
    @{ Cell = 'E25'; Value = '  +3.66%  ' },
    @{ Cell = 'D26'; Value = '2.409' },
    @{ Cell = 'E26'; Value = '  -6.35%  ' },
    @{ Cell = 'D27'; Value = '148.67' },
    @{ Cell = 'E27'; Value = '  -1.37%  ' },
    @{ Cell = 'D28'; Value = '18.48' },
    @{ Cell = 'E28'; Value = '  -4.29%  ' },
    @{ Cell = 'D29'; Value = '4.883' },
    @{ Cell = 'E29'; Value = '  -1.27%  ' },
    @{ Cell = 'D30'; Value = '1.726.23' },
    @{ Cell = 'E30'; Value = '  -0.09%  ' },
    @{ Cell = 'D31'; Value = '118.36' },
    @{ Cell = 'E31'; Value = '  -3.91%  ' },
    @{ Cell = 'D32'; Value = '0.9892' },
    @{ Cell = 'E32'; Value = '  -8.56%  ' },
    @{ Cell = 'D33'; Value = '5.789' },
    @{ Cell = 'E33'; Value = '  -2.38%  ' },
    @{ Cell = 'D34'; Value = '0.08285' },
    @{ Cell = 'E34'; Value = '  +0.74%  ' },
    @{ Cell = 'D35'; Value = '1.613' },
    @{ Cell = 'E35'; Value = '  -16.03%  ' },
    @{ Cell = 'D36'; Value = '8.982' },
    @{ Cell = 'E36'; Value = '  -4.86%  ' },
    @{ Cell = 'D37'; Value = '0.02253' },
    @{ Cell = 'E37'; Value = '  -4.62%  ' },
    @{ Cell = 'D38'; Value = '0.06057' },
    @{ Cell = 'E38'; Value = '  -4.72%  ' },
    @{ Cell = 'D39'; Value = '5.086' },
    @{ Cell = 'E39'; Value = '  -5.24%  ' },
    @{ Cell = 'D40'; Value = '1.210' },
    @{ Cell = 'E40'; Value = '  -2.63%  ' },
    @{ Cell = 'D41'; Value = '0.2035' },
    @{ Cell = 'E41'; Value = '  -5.94%  ' },
    @{ Cell = 'E42'; Value = '  +0.05%  ' },
    @{ Cell = 'D43'; Value = '10.61' },
    @{ Cell = 'E43'; Value = '  -3.98%  ' },
    @{ Cell = 'D44'; Value = '0.5787' },
    @{ Cell = 'E44'; Value = '  -4.96%  ' },
    @{ Cell = 'D45'; Value = '3.745' },
    @{ Cell = 'E45'; Value = '  -0.25%  ' },
    @{ Cell = 'D46'; Value = '12.91' },
    @{ Cell = 'E46'; Value = '  -7.16%  ' },
    @{ Cell = 'D47'; Value = '0.5557' },
    @{ Cell = 'E47'; Value = '  -6.00%  ' },
    @{ Cell = 'D48'; Value = '117.97' },
    @{ Cell = 'E48'; Value = '  -4.38%  ' },
    @{ Cell = 'D49'; Value = '1.886' },
    @{ Cell = 'E49'; Value = '  -4.58%  ' },
    @{ Cell = 'D50'; Value = '1.130' },
    @{ Cell = 'E50'; Value = '  -4.15%  ' },
    @{ Cell = 'D51'; Value = '0.06812' },
    @{ Cell = 'E51'; Value = '  -3.78%  ' }
)

foreach ($u in $updates) {
    $cell = $ws.Range($u.Cell)
    # Force text interpretation so values like "1.003" or "290.17" are not
    # silently coerced into numbers by Excel, then strip the format change
    # back off so the cell's style stays exactly as it was (no explicit
    # style index), matching the original inline-string cells.
    $cell.NumberFormat = "@"
    $cell.Value = $u.Value
    $cell.ClearFormats()
}
